$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("P2").NumberFormat = "@"
$ws.Range("P2").Value = "320018253468"
$ws.Range("P2").Style = "Normal"
$ws.Range("P3").NumberFormat = "@"
$ws.Range("P3").Value = "320018253479"
$ws.Range("P3").Style = "Normal"
